$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12534
$ws1.Range("C4").Value = "苏州·ICOMEACG展后回馈x星渡咖啡（取消）"
$ws1.Range("F4").Value = 2035
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F8").Value = 12502
$ws1.Range("F9").Value = 3080
$ws1.Range("F10").Value = 532
$ws1.Range("F14").Value = 22
$ws1.Range("F16").Value = 646
$ws1.Range("F22").Value = 37

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12534
$ws4.Range("C4").Value = "苏州·ICOMEACG展后回馈x星渡咖啡（取消）"
$ws4.Range("F4").Value = 2035
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F9").Value = 12502
$ws4.Range("F10").Value = 3080
$ws4.Range("F11").Value = 532
$ws4.Range("F15").Value = 22
$ws4.Range("F17").Value = 646
$ws4.Range("F24").Value = 37
